# [MOSIP-14369] Fix: boolean values
# Replace the TRUE() formulas in column D (is_active) with the literal text "TRUE"
# so the value is stored as a string instead of a computed boolean number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("D2:D13")

# Replace the TRUE() boolean formulas with the literal text string "TRUE".
# A formula that evaluates to the text "TRUE" is written first, then the
# range is pasted back onto itself as values only so the formula is
# discarded and a plain text cell (shared string) remains.
$range.Formula = '="TRUE"'
$range.Copy() | Out-Null
$range.PasteSpecial(-4163) | Out-Null        # xlPasteValues
$excel.CutCopyMode = 0

# Update the active selection to mirror the edited range
$range.Select() | Out-Null
